# Updates the "To do.docx" to-do list:
#  - adds several new bullet items before the "Bugeja:" heading
#  - adds a new bug entry (and a trailing blank-ish paragraph) after "Bugeja:"

function New-FlatOpcBodyXml($bodyInnerXml) {
    # Minimal Flat-OPC wrapper so Range.InsertXML can graft a fully-formed
    # <w:p> (with whatever run children we like, e.g. <w:tab/>) into the
    # document body at the target range.
    return "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
           "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
           "<pkg:xmlData>" +
           "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
           "<w:body>$bodyInnerXml</w:body>" +
           "</w:document>" +
           "</pkg:xmlData></pkg:part></pkg:package>"
}

function Add-ParagraphAfter($para, $text, [bool]$leadingTab) {
    # Inserts a brand new paragraph right after $para, containing $text
    # (optionally preceded by a real <w:tab/> run child), and returns the
    # new Paragraph object so callers can keep chaining.
    $para.Range.InsertParagraphAfter()
    $newPara = $para.Next()

    $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $space = ""
    if ($text -ne $text.Trim()) {
        $space = " xml:space='preserve'"
    }
    if ($leadingTab) {
        $runXml = "<w:r><w:tab/><w:t$space>$escaped</w:t></w:r>"
    } else {
        $runXml = "<w:r><w:t$space>$escaped</w:t></w:r>"
    }
    $xml = New-FlatOpcBodyXml("<w:p>$runXml</w:p>")
    $newPara.Range.InsertXML($xml)

    return $para.Next()
}

$d = $word.ActiveDocument

# Anchor on the "-Tuli tarttumaan puusta toiseen" paragraph (just before
# "Bugeja:") and append the six new to-do lines after it.
$cur = $d.Paragraphs.Item(4)

$cur = Add-ParagraphAfter $cur "-Animaatioita lisää" $false
$cur = Add-ParagraphAfter $cur "-Decalit poolista(ei siis instantiatella, vaan siten että niitä ladataan tietty määrä johonkin kauas ja sitten siirretään tarvittavaan paikkaan telaketjumalliin):" $false
$cur = Add-ParagraphAfter $cur "-Luodinreiät maahan, sekä seiniin ja muihin esineisiin" $true
$cur = Add-ParagraphAfter $cur "-Verinen luodinreikä hahmoihin, jos onnistuu kohtuullisella vaivalla, tuskin onnistuu" $true
$cur = Add-ParagraphAfter $cur "-Veriroiske(bloodsplash) aina osuman jälkeen maahan" $true
$cur = Add-ParagraphAfter $cur "-Räjähdyskraaterit maahan räjähdyksestä" $true

# The trailing empty paragraph (after "Bugeja:") becomes the new bug entry,
# and a fresh, single-space paragraph is appended after it.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.Text = "-Grafiikka-asetukset vaikuttavat hyppykorkeuteen... :-D"

$last.Range.InsertParagraphAfter()
$newLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$newLast.Range.Text = " "

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
